$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.379.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.64%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.397.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "405.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.89"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +15.87%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +8.02%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.49%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +10.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.36"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.59%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.29%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.949.31"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.55"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.19%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.72"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.399.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.52"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +9.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.307.81"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000134"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +18.29%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "82.92"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +13.69%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "308.51"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.79%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.87%  "

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.75"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +10.38%  "

# Row 27
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.54"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +15.18%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.79%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.85%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.91%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.43%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.82"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.54%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.72"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +9.43%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.22%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.52%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.33"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.25%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.12%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.51%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.25%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.49"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.12%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.20%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.05"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.33%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.82"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.151.28"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.19%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.729.08"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.27%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.47%  "
